$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q1" sheet by copying the "2021-Q4" sheet's
#     layout/formatting, positioned right before the "总计" sheet. ---
$src = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")
$src.Copy($totalSheet, $null)

$ws = $wb.Worksheets.Item("2021-Q4 (2)")
$ws.Name = "2022-Q1"

# Header row
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Data rows. Fund codes (B) and the metric columns (D:G) are stored as text
# in the source data (leading zeros on fund codes, fixed-decimal formatting
# on the metrics), so they're entered with a leading apostrophe to keep
# Excel from re-parsing them as numbers.
$data = @(
    @(0, "002692", "富国创新科技混合A", "40.49", "92.07", "5.07", "2.0528", 5),
    @(1, "100060", "富国高新技术产业混合", "37.64", "92.18", "4.11", "1.5470", 7),
    @(2, "501077", "富国科创主题 3 年封闭运作灵活配置混 合型", "19.40", "99.17", "4.90", "0.9506", 6),
    @(3, "320005", "诺安价值增长混合", "10.72", "91.23", "6.88", "0.7375", 4),
    @(4, "007345", "富国科技创新灵活配置混合", "12.61", "92.32", "4.91", "0.6192", 5),
    @(5, "002291", "诺安安鑫灵活配置混合", "2.19", "81.55", "7.49", "0.1640", 3),
    @(6, "011120", "富国创新科技混合C", "0.90", "92.07", "5.07", "0.0456", 5)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = "'" + $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = "'" + $rowData[3]
    $ws.Cells.Item($row, 5).Value = "'" + $rowData[4]
    $ws.Cells.Item($row, 6).Value = "'" + $rowData[5]
    $ws.Cells.Item($row, 7).Value = "'" + $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
}

# The leading apostrophe marks the cell "number stored as text" (quote-
# prefix) in its style; clear that back to the sheet's plain default style
# on the text columns, same as the source data's cells carry no special
# style. Column A keeps its own centered/bordered style ("s=2").
$ws.Range("Z1").Copy()
$ws.Range("B2:B8").PasteSpecial(-4122)
$ws.Range("D2:G8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("A3:A8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 2: insert a new top data-row in "总计" sheet for 2022-Q1, shifting
#     the existing rows down by one whole row. ---
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 6.12

# The freshly inserted row doesn't carry the original column formatting
# (A = bordered/centered "s=2", B:D = plain/no style) -- restore it from the
# row directly below, which still has it.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Range("B3:D3").Copy()
$total.Range("B2:D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Renumber the index column (A) for the rows that were pushed down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
